$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '281.75'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '20.67'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06143'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.576'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.563'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.495'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8164'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1634'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08321'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03542'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03184'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09133'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.707'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04702'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006520'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006160'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001501'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.779'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.322'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3359'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04672'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007196'

$ws.Range("B42").Value = 'BKEXToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1100'

$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'

$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003513'

$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01102'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006627'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00001901'
